$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $s = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $s
}

$ws.Range("D2").Value = '22.452.14'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.568.28'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  -0.19%  '
Set-TextValue $ws "D5" '1.000'
$ws.Range("E5").Value = '  -0.24%  '
Set-TextValue $ws "D6" '289.01'
$ws.Range("E6").Value = '  -0.36%  '
Set-TextValue $ws "D7" '0.3686'
$ws.Range("E7").Value = '  -1.74%  '
Set-TextValue $ws "D8" '49.81'
$ws.Range("E8").Value = '  -0.06%  '
Set-TextValue $ws "D9" '0.3382'
$ws.Range("E9").Value = '  +0.13%  '
Set-TextValue $ws "D10" '1.135'
$ws.Range("E10").Value = '  +0.07%  '
Set-TextValue $ws "D11" '0.07483'
$ws.Range("E11").Value = '  -0.23%  '
Set-TextValue $ws "D12" '1.001'
$ws.Range("E12").Value = '  -0.23%  '
Set-TextValue $ws "D13" '21.02'
$ws.Range("E13").Value = '  -2.12%  '
Set-TextValue $ws "D14" '5.974'
$ws.Range("E14").Value = '  +0.29%  '
Set-TextValue $ws "D15" '6.939'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").Value = '1.565.45'
$ws.Range("E16").Value = '  -0.49%  '
Set-TextValue $ws "D17" '0.00001109'
$ws.Range("E17").Value = '  -0.56%  '
Set-TextValue $ws "D18" '90.15'
$ws.Range("E18").Value = '  +0.30%  '
Set-TextValue $ws "D19" '0.06757'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  -0.22%  '
Set-TextValue $ws "D21" '6.366'
$ws.Range("E21").Value = '  +2.57%  '
Set-TextValue $ws "D22" '16.13'
$ws.Range("E22").Value = '  -0.75%  '
Set-TextValue $ws "D23" '12.03'
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").Value = '22.441.53'
$ws.Range("E24").Value = '  +0.18%  '
Set-TextValue $ws "D25" '2.388'
$ws.Range("E25").Value = '  +0.86%  '
Set-TextValue $ws "D26" '2.588'
$ws.Range("E26").Value = '  +2.15%  '
Set-TextValue $ws "D27" '19.76'
$ws.Range("E27").Value = '  -2.01%  '
Set-TextValue $ws "D28" '149.70'
$ws.Range("E28").Value = '  +1.36%  '
Set-TextValue $ws "D29" '5.029'
$ws.Range("E29").Value = '  +0.89%  '
Set-TextValue $ws "D30" '124.12'
$ws.Range("E30").Value = '  -0.75%  '
$ws.Range("D31").Value = '1.740.07'
$ws.Range("E31").Value = '  -0.50%  '
Set-TextValue $ws "D32" '1.059'
$ws.Range("E32").Value = '  +5.60%  '
Set-TextValue $ws "D33" '6.142'
$ws.Range("E33").Value = '  +3.01%  '
Set-TextValue $ws "D34" '2.010'
$ws.Range("E34").Value = '  +2.41%  '
Set-TextValue $ws "D35" '9.650'
$ws.Range("E35").Value = '  -0.92%  '
Set-TextValue $ws "D36" '0.08311'
$ws.Range("E36").Value = '  -1.08%  '
Set-TextValue $ws "D37" '0.02432'
$ws.Range("E37").Value = '  -0.50%  '
Set-TextValue $ws "D38" '1.331'
$ws.Range("E38").Value = '  -3.49%  '
Set-TextValue $ws "D39" '0.2234'
$ws.Range("E39").Value = '  -0.82%  '
Set-TextValue $ws "D40" '0.06382'
$ws.Range("E40").Value = '  -2.09%  '
Set-TextValue $ws "D41" '5.353'
$ws.Range("E41").Value = '  -1.14%  '
Set-TextValue $ws "D42" '11.18'
$ws.Range("E42").Value = '  -1.07%  '
Set-TextValue $ws "D43" '0.6121'
$ws.Range("E43").Value = '  -1.74%  '
Set-TextValue $ws "D44" '0.9999'
$ws.Range("E44").Value = '  -0.15%  '
Set-TextValue $ws "D45" '13.91'
$ws.Range("E45").Value = '  -1.13%  '
Set-TextValue $ws "D46" '3.768'
$ws.Range("E46").Value = '  -1.19%  '
Set-TextValue $ws "D47" '0.5761'
$ws.Range("E47").Value = '  -0.58%  '
Set-TextValue $ws "D48" '2.025'
$ws.Range("E48").Value = '  -2.72%  '
Set-TextValue $ws "D49" '125.59'
$ws.Range("E49").Value = '  -0.95%  '
Set-TextValue $ws "D50" '1.233'
$ws.Range("E50").Value = '  +0.53%  '
Set-TextValue $ws "D51" '0.07312'
$ws.Range("E51").Value = '  +0.07%  '
